$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("en")
$ws2 = $wb.Worksheets.Item("de")

# ----------------------------------------------------------------------------
# Pretty checkbox library: new login/registration resource strings.
# Step 1 - write all new cell VALUES first, in the exact order the original
#          authoring happened, so new entries land in xl/sharedStrings.xml in
#          the same sequence as the target workbook.
# ----------------------------------------------------------------------------

$ws1.Range("A109").Value = "CreateAccount"
$ws1.Range("B109").Value = "Create a new account."
$ws1.Range("A110").Value = "Password"
$ws2.Range("B110").Value = "Passwort"
$ws2.Range("B109").Value = "Ein neues Konto eröffnen."
$ws1.Range("B111").Value = "The {0} must be at least {2} and at max {1} characters long."
$ws1.Range("A111").Value = "PasswordErrorLength"
$ws2.Range("B111").Value = "Das {0} muss mindestens {2} und höchstens {1} Zeichen lang sein."
$ws1.Range("A112").Value = "ConfirmPassword"
$ws2.Range("B112").Value = "Passwort bestätigen"
$ws1.Range("B112").Value = "Confirm password"
$ws1.Range("A113").Value = "CompareFailed"
$ws1.Range("B113").Value = "The password and confirmation password do not match."
$ws2.Range("B113").Value = "Die Passwörter stimmen nicht überein"
$ws1.Range("A114").Value = "LoginTitle"
$ws1.Range("B114").Value = "Use a local account to log in."
$ws2.Range("B114").Value = "Melden Sie sich mit einem Konto an."
$ws1.Range("B115").Value = "Remember me?"
$ws2.Range("B115").Value = "Eingeloggt bleiben"
$ws1.Range("A115").Value = "RememberMe"
$ws1.Range("B110").Value = "Password"
$ws2.Range("A109").Value = "CreateAccount"
$ws2.Range("A110").Value = "Password"
$ws2.Range("A111").Value = "PasswordErrorLength"
$ws2.Range("A112").Value = "ConfirmPassword"
$ws2.Range("A113").Value = "CompareFailed"
$ws2.Range("A114").Value = "LoginTitle"
$ws2.Range("A115").Value = "RememberMe"

# ----------------------------------------------------------------------------
# Step 2 - apply the "wrap text + vertically centered" look (same visual style
#          used throughout the sheet, cellXfs index 1) to the cells that need
#          it, by copying the format from an already-styled template cell.
#          This reuses the existing style instead of registering new ones.
# ----------------------------------------------------------------------------

$ws1.Range("A2").Copy() | Out-Null
$ws1.Range("A109").PasteSpecial(-4122) | Out-Null
$ws1.Range("A110").PasteSpecial(-4122) | Out-Null
$ws1.Range("A111").PasteSpecial(-4122) | Out-Null
$ws1.Range("B112").PasteSpecial(-4122) | Out-Null
$ws1.Range("A113").PasteSpecial(-4122) | Out-Null
$ws1.Range("A114").PasteSpecial(-4122) | Out-Null
$ws1.Range("B115").PasteSpecial(-4122) | Out-Null
$ws1.Range("A115").PasteSpecial(-4122) | Out-Null
$ws1.Range("B110").PasteSpecial(-4122) | Out-Null

$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("B110").PasteSpecial(-4122) | Out-Null
$ws2.Range("B109").PasteSpecial(-4122) | Out-Null
$ws2.Range("B112").PasteSpecial(-4122) | Out-Null
$ws2.Range("A109").PasteSpecial(-4122) | Out-Null
$ws2.Range("A110").PasteSpecial(-4122) | Out-Null
$ws2.Range("A111").PasteSpecial(-4122) | Out-Null
$ws2.Range("A113").PasteSpecial(-4122) | Out-Null
$ws2.Range("A114").PasteSpecial(-4122) | Out-Null
$ws2.Range("A115").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ----------------------------------------------------------------------------
# Step 3 - restore view state: active sheet, scroll position, selection.
# ----------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 94
$ws1.Range("A115:B115").Select() | Out-Null

$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 94
$ws2.Range("A115:B115").Select() | Out-Null

